# The author retyped/confirmed the final period of the sentence that ends
# "... les problèmes dans lesquels ils sont utilisés." — which is exactly
# the kind of edit that makes Word drop its "last edit" (_GoBack) bookmark
# right before that period, removing it from wherever it used to sit
# (in this document: the empty paragraph at the very end).
#
# Re-creating a bookmark named "_GoBack" automatically replaces/moves any
# existing bookmark of the same name, so we only need to add the new one;
# Word takes care of removing the stale one for us.

$d = $word.ActiveDocument

$matchRange = $d.Content
$found = $matchRange.Find.Execute("utilisés.", $true, $false, $false, $false, `
                                   $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not locate the text 'utilisés.' in the document."
}

# Collapse to the zero-length point right before the trailing "." so the
# new bookmark lands between "s" and "." — this is what causes Word to
# split that run into "s" + "." when it next serializes the paragraph.
$periodStart = $matchRange.End - 1
$goBackPoint = $d.Range($periodStart, $periodStart)

$d.Bookmarks.Add("_GoBack", $goBackPoint)
